# BOM-byvalue.xlsx edit:
# Replace the single H-bridge logic translator line (TXB0108, QFN20, qty 3, U17)
# with a new DQS variant line (TXB0108DQSR, DQS_R-PUSON-N20, qty 2, for U12/U19),
# leaving one TXB0108 QFN unit (qty 1) on the original U17 line.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 92 (HMC5883LSMD / U13 line).
# This pushes the existing rows 92-102 down to 93-103, and Excel automatically
# adjusts the dependent formulas/ranges (K2's SUM, the shared-formula ranges,
# and the sheet dimension).
$ws.Rows.Item(92).Insert()

# Populate the newly inserted row 92 with the DQS part line.
# Set column D (Package) before column A (Part/designator) so the shared
# string table is built in the same order as the source edit.
$ws.Cells.Item(92, 4).Value = "DQS_R-PUSON-N20"
$ws.Cells.Item(92, 1).Value = "U12, U19"
$ws.Cells.Item(92, 2).Value = "TXB0108"
$ws.Cells.Item(92, 3).Value = "TXB0108"
$ws.Cells.Item(92, 5).Value = "ME"
$ws.Cells.Item(92, 6).Value = "595-TXB0108DQSR"
$ws.Cells.Item(92, 7).Value = "TXB0108DQSR"
$ws.Cells.Item(92, 8).Value = 2
$ws.Cells.Item(92, 9).Value = 2.76
$ws.Cells.Item(92, 10).Formula = "=H92*I92"

# The original TXB0108 / QFN20 / U17 line (now shifted to row 97) only keeps
# one unit instead of three, since two of the three are replaced by the new
# DQS variant line above.
$ws.Cells.Item(97, 8).Value = 1

# Restore the view state: scrolled so row 64 is the top-left visible row,
# with I97 as the active selection.
$ws.Activate()
$ws.Range("I97").Select()
$excel.ActiveWindow.ScrollRow = 64
$excel.ActiveWindow.ScrollColumn = 1
